# Applies the weekly Fruta/Hortaliza price-sheet refresh for the
# "Feria Lagunitas de Puerto Montt - Granada" subset: each data row (2-39)
# gets its date, quality grade, volume, min/max/avg price, unit label,
# $/Kg price and Kg/unit fields updated to the new week's figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44301
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 17500
$ws.Range("O2").Value = 17500
$ws.Range("P2").Value = 17500
$ws.Range("S2").Value = 1167

# Row 3
$ws.Range("D3").Value = 44301
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 14000
$ws.Range("P3").Value = 14500
$ws.Range("S3").Value = 967

# Row 4
$ws.Range("D4").Value = 44334
$ws.Range("L4").Value = "Primera"
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 15500
$ws.Range("S4").Value = 1033

# Row 5
$ws.Range("D5").Value = 44334
$ws.Range("L5").Value = "Segunda"
$ws.Range("N5").Value = 14500
$ws.Range("O5").Value = 14500
$ws.Range("P5").Value = 14500
$ws.Range("Q5").Value = "$/caja 15 kilos empedrada"
$ws.Range("S5").Value = 967
$ws.Range("T5").Value = 15

# Row 6
$ws.Range("D6").Value = 44330
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 17000
$ws.Range("P6").Value = 17000
$ws.Range("Q6").Value = "$/caja 15 kilos empedrada"
$ws.Range("S6").Value = 1133
$ws.Range("T6").Value = 15

# Row 7
$ws.Range("D7").Value = 44330
$ws.Range("M7").Value = 200
$ws.Range("O7").Value = 14500
$ws.Range("P7").Value = 14250
$ws.Range("Q7").Value = "$/caja 15 kilos empedrada"
$ws.Range("S7").Value = 950
$ws.Range("T7").Value = 15

# Row 8
$ws.Range("D8").Value = 44351
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("S8").Value = 1000

# Row 9
$ws.Range("D9").Value = 44351
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 13000
$ws.Range("O9").Value = 13500
$ws.Range("P9").Value = 13250
$ws.Range("S9").Value = 883

# Row 10
$ws.Range("D10").Value = 44302

# Row 11
$ws.Range("D11").Value = 44302
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 14500
$ws.Range("S11").Value = 967

# Row 12
$ws.Range("D12").Value = 44348
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("S12").Value = 1000

# Row 13
$ws.Range("D13").Value = 44348
$ws.Range("M13").Value = 200

# Row 14
$ws.Range("D14").Value = 44309
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 17500
$ws.Range("O14").Value = 17500
$ws.Range("P14").Value = 17500
$ws.Range("S14").Value = 1167

# Row 15
$ws.Range("D15").Value = 44309
$ws.Range("L15").Value = "Segunda"
$ws.Range("M15").Value = 200
$ws.Range("O15").Value = 14500
$ws.Range("P15").Value = 14250
$ws.Range("S15").Value = 950

# Row 16
$ws.Range("D16").Value = 44323
$ws.Range("N16").Value = 17000
$ws.Range("O16").Value = 17000
$ws.Range("P16").Value = 17000
$ws.Range("S16").Value = 1133

# Row 17
$ws.Range("D17").Value = 44323
$ws.Range("M17").Value = 100
$ws.Range("O17").Value = 14000
$ws.Range("P17").Value = 14000
$ws.Range("S17").Value = 933

# Row 18
$ws.Range("D18").Value = 44295
$ws.Range("L18").Value = "Primera"

# Row 19
$ws.Range("D19").Value = 44327
$ws.Range("N19").Value = 17000
$ws.Range("O19").Value = 17000
$ws.Range("P19").Value = 17000
$ws.Range("S19").Value = 1133

# Row 20
$ws.Range("D20").Value = 44327

# Row 21
$ws.Range("D21").Value = 44305
$ws.Range("M21").Value = 60
$ws.Range("N21").Value = 17500
$ws.Range("O21").Value = 17500
$ws.Range("P21").Value = 17500
$ws.Range("S21").Value = 1167

# Row 22
$ws.Range("D22").Value = 44305
$ws.Range("M22").Value = 120
$ws.Range("O22").Value = 15000
$ws.Range("P22").Value = 14500
$ws.Range("S22").Value = 967

# Row 23
$ws.Range("D23").Value = 44293
$ws.Range("M23").Value = 60

# Row 24
$ws.Range("D24").Value = 44336
$ws.Range("N24").Value = 17000
$ws.Range("O24").Value = 17000
$ws.Range("P24").Value = 17000
$ws.Range("S24").Value = 1133

# Row 25
$ws.Range("D25").Value = 44336
$ws.Range("O25").Value = 14500
$ws.Range("P25").Value = 14250
$ws.Range("S25").Value = 950

# Row 26
$ws.Range("D26").Value = 44292
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 160
$ws.Range("N26").Value = 14000
$ws.Range("O26").Value = 15000
$ws.Range("P26").Value = 14500
$ws.Range("S26").Value = 967

# Row 27
$ws.Range("D27").Value = 44306
$ws.Range("L27").Value = "Primera"
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 17500
$ws.Range("O27").Value = 17500
$ws.Range("P27").Value = 17500
$ws.Range("S27").Value = 1167

# Row 28
$ws.Range("D28").Value = 44306
$ws.Range("L28").Value = "Segunda"
$ws.Range("M28").Value = 200
$ws.Range("N28").Value = 14000
$ws.Range("O28").Value = 14500
$ws.Range("P28").Value = 14250
$ws.Range("S28").Value = 950

# Row 29
$ws.Range("D29").Value = 44285
$ws.Range("M29").Value = 160
$ws.Range("N29").Value = 15000
$ws.Range("O29").Value = 16000
$ws.Range("P29").Value = 15500
$ws.Range("S29").Value = 1033

# Row 30
$ws.Range("D30").Value = 44313
$ws.Range("L30").Value = "Especial"
$ws.Range("M30").Value = 100
$ws.Range("Q30").Value = "$/caja 14 kilos empedrada"
$ws.Range("S30").Value = 1250
$ws.Range("T30").Value = 14

# Row 31
$ws.Range("D31").Value = 44313
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 100
$ws.Range("N31").Value = 16000
$ws.Range("O31").Value = 16000
$ws.Range("P31").Value = 16000
$ws.Range("Q31").Value = "$/caja 14 kilos empedrada"
$ws.Range("S31").Value = 1143
$ws.Range("T31").Value = 14

# Row 32
$ws.Range("D32").Value = 44313
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 80
$ws.Range("N32").Value = 14000
$ws.Range("O32").Value = 14000
$ws.Range("P32").Value = 14000
$ws.Range("Q32").Value = "$/caja 14 kilos empedrada"
$ws.Range("T32").Value = 14

# Row 33
$ws.Range("D33").Value = 44298
$ws.Range("M33").Value = 80
$ws.Range("N33").Value = 14000
$ws.Range("O33").Value = 15000
$ws.Range("P33").Value = 14500
$ws.Range("S33").Value = 967

# Row 34
$ws.Range("D34").Value = 44299
$ws.Range("M34").Value = 60
$ws.Range("N34").Value = 17500
$ws.Range("O34").Value = 17500
$ws.Range("P34").Value = 17500
$ws.Range("S34").Value = 1167

# Row 35
$ws.Range("D35").Value = 44299
$ws.Range("M35").Value = 120
$ws.Range("O35").Value = 15000
$ws.Range("P35").Value = 14500
$ws.Range("S35").Value = 967

# Row 36
$ws.Range("D36").Value = 44316
$ws.Range("N36").Value = 17500
$ws.Range("O36").Value = 17500
$ws.Range("P36").Value = 17500
$ws.Range("S36").Value = 1167

# Row 37
$ws.Range("D37").Value = 44316

# Row 38
$ws.Range("D38").Value = 44344
$ws.Range("N38").Value = 16000
$ws.Range("O38").Value = 16000
$ws.Range("P38").Value = 16000
$ws.Range("S38").Value = 1067

# Row 39
$ws.Range("D39").Value = 44344
$ws.Range("M39").Value = 120
$ws.Range("N39").Value = 13000
$ws.Range("O39").Value = 13500
$ws.Range("P39").Value = 13250
$ws.Range("S39").Value = 883
